$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '37.957.86'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  -0.02%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.037.21'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  -0.68%  '

# Row 4
$ws.Range('E4').Value = '  -0.12%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.34'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.56%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.611'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.56%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '60.61'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +3.45%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.380'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.40%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0822'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +1.30%  '

# Row 11
$ws.Range('E11').Value = '  +0.01%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.337.76'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -0.75%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.56'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.24%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.43'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +2.71%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.764'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +1.67%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.18'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.73%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.052.63'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.35%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '37.906.79'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.01%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.73'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.07%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.92'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -6.55%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.0₃0827'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.06%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '224.21'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -0.09%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.999'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -0.08%  '

# Row 24
$ws.Range('E24').Value = '  +0.04%  '

# Row 25
$ws.Range('E25').Value = '  -0.15%  '

# Row 26
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '9.32'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.65%  '

# Row 27
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '166.93'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.25%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.131'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.75%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.86'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.77%  '

# Row 30
$ws.Range('E30').Value = '  -3.35%  '

# Row 31
$ws.Range('E31').Value = '  +0.83%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.29'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +10.91%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.42'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.30%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0607'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +0.12%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.52'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.63%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.40'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +5.73%  '

# Row 37
$ws.Range('E37').Value = '  -1.76%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.29'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.78%  '

# Row 39
$ws.Range('E39').Value = '  +0.08%  '

# Row 40
$ws.Range('B40').Value = 'InjectiveProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '17.71'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +7.27%  '

# Row 41
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.538.56'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.28%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0219'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.81%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '96.41'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.45%  '

# Row 44
$ws.Range('E44').Value = '  -3.42%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0914'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -0.45%  '

# Row 46
$ws.Range('E46').Value = '  -2.37%  '

# Row 47
$ws.Range('E47').Value = '  -1.56%  '

# Row 48
$ws.Range('E48').Value = '  -0.43%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.97'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.08%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.09'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.03%  '

# Row 51
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.226.11'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.75%  '
